# This workbook is a weekly-updated consolidated price series for
# "Femacal de La Calera - Cebolla" (onion prices). A new reporting week
# (fecha serial 44516) is inserted as 3 new data rows at the top of the
# data block (rows 765-767), pushing the previously existing rows
# 765-814 down to 768-817 and extending the sheet dimension to A1:R817.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows before row 765; this shifts everything at/after
# row 765 down by three rows (the date style on column D carries down with
# Insert, matching how Excel's Rows.Insert behaves by default).
$ws.Rows("765:767").Insert()

# New row 765 (fecha = 44516)
$ws.Range("A765").Value = 3
$ws.Range("B765").Value = 'Femacal de La Calera'
$ws.Range("C765").Value = 'Coquimbo'
$ws.Range("D765").Value = 44516
$ws.Range("E765").Value = 5
$ws.Range("F765").Value = 100112004
$ws.Range("G765").Value = 'Cebolla'
$ws.Range("H765").Value = 'Morada(o)'
$ws.Range("I765").Value = '1a nueva(o)'
$ws.Range("J765").Value = 105
$ws.Range("K765").Value = 5500
$ws.Range("L765").Value = 6000
$ws.Range("M765").Value = 5738
$ws.Range("N765").Value = '$/malla 18 kilos'
$ws.Range("O765").Value = 'Región de Arica y Parinacota'
$ws.Range("P765").Value = 319
$ws.Range("Q765").Value = 18
$ws.Range("R765").Value = 'Hortaliza'

# New row 766 (fecha = 44516)
$ws.Range("A766").Value = 3
$ws.Range("B766").Value = 'Femacal de La Calera'
$ws.Range("C766").Value = 'Coquimbo'
$ws.Range("D766").Value = 44516
$ws.Range("E766").Value = 5
$ws.Range("F766").Value = 100112004
$ws.Range("G766").Value = 'Cebolla'
$ws.Range("H766").Value = 'Sin especificar'
$ws.Range("I766").Value = '1a nueva(o)'
$ws.Range("J766").Value = 125
$ws.Range("K766").Value = 4500
$ws.Range("L766").Value = 4700
$ws.Range("M766").Value = 4596
$ws.Range("N766").Value = '$/malla 18 kilos'
$ws.Range("O766").Value = 'Provincia de Quillota'
$ws.Range("P766").Value = 255
$ws.Range("Q766").Value = 18
$ws.Range("R766").Value = 'Hortaliza'

# New row 767 (fecha = 44516)
$ws.Range("A767").Value = 3
$ws.Range("B767").Value = 'Femacal de La Calera'
$ws.Range("C767").Value = 'Coquimbo'
$ws.Range("D767").Value = 44516
$ws.Range("E767").Value = 5
$ws.Range("F767").Value = 100112004
$ws.Range("G767").Value = 'Cebolla'
$ws.Range("H767").Value = 'Sin especificar'
$ws.Range("I767").Value = '1a nueva(o)'
$ws.Range("J767").Value = 120
$ws.Range("K767").Value = 4000
$ws.Range("L767").Value = 4000
$ws.Range("M767").Value = 4000
$ws.Range("N767").Value = '$/malla 18 kilos'
$ws.Range("O767").Value = 'Región de Arica y Parinacota'
$ws.Range("P767").Value = 222
$ws.Range("Q767").Value = 18
$ws.Range("R767").Value = 'Hortaliza'
